# Applies updated "want to go" counts (column F) from the site regeneration commit
# "Update gh-pages to output generated at 456a3b4"
$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml) - column F ("想去人数") updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 1141
$ws1.Range("F9").Value = 120
$ws1.Range("F14").Value = 830
$ws1.Range("F15").Value = 844
$ws1.Range("F20").Value = 721
$ws1.Range("F21").Value = 1725
$ws1.Range("F22").Value = 2596
$ws1.Range("F23").Value = 738
$ws1.Range("F25").Value = 2005
$ws1.Range("F26").Value = 461
$ws1.Range("F27").Value = 2892
$ws1.Range("F29").Value = 85
$ws1.Range("F30").Value = 708
$ws1.Range("F31").Value = 137
$ws1.Range("F32").Value = 114
$ws1.Range("F34").Value = 1008
$ws1.Range("F35").Value = 1727
$ws1.Range("F36").Value = 361
$ws1.Range("F38").Value = 543

# Sheet "全部类型" (sheet4.xml) - column F ("想去人数") updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 1141
$ws4.Range("F10").Value = 120
$ws4.Range("F14").Value = 830
$ws4.Range("F15").Value = 845
$ws4.Range("F21").Value = 721
$ws4.Range("F22").Value = 1725
$ws4.Range("F23").Value = 2596
$ws4.Range("F24").Value = 738
$ws4.Range("F28").Value = 2892
$ws4.Range("F34").Value = 85
$ws4.Range("F36").Value = 708
$ws4.Range("F37").Value = 137
$ws4.Range("F38").Value = 114
$ws4.Range("F40").Value = 1008
$ws4.Range("F41").Value = 1727
$ws4.Range("F43").Value = 361
$ws4.Range("F44").Value = 543
